# "Generate Report for Handoff" - the localization status flips from
# "In Translation" to "Ready for handoff" and the associated timestamps
# advance, for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2/F2 = zh-cn / de-de Status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-18 12:59:11"

# --- zh-cn sheet -------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-18 12:59:00"

# --- de-de sheet -------------------------------------------------------
# C2 = Status (H2 keeps referencing the same "Latest HO Xliff Generate
# Date" text as Overview!G2, so it updates to the same new timestamp)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-18 12:59:11"

# --- Column widths -------------------------------------------------------
# The Status columns grow because "Ready for handoff" renders wider than
# "In Translation"; set the new widened width directly.
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
